$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 301
$ws1.Range("F12").Value = 184
$ws1.Range("F26").Value = 6061
$ws1.Range("F27").Value = 75
$ws1.Range("F31").Value = 14775
$ws1.Range("F34").Value = 111
$ws1.Range("F36").Value = 10816

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 301
$ws4.Range("F12").Value = 184
$ws4.Range("F29").Value = 6061
$ws4.Range("F30").Value = 75
$ws4.Range("F34").Value = 14775
$ws4.Range("F37").Value = 111
$ws4.Range("F39").Value = 10816
